# Updated cryptos list on Fri Oct 27 05:49:55 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Price($row, $value) {
    # Prefix with an apostrophe so values that look numeric (e.g. "224.74")
    # are stored as literal text, matching the source inlineStr cells.
    $ws.Cells.Item($row, 4).Value = "'" + $value
}

function Set-Volume($row, $value) {
    $ws.Cells.Item($row, 5).Value = $value
}

# Row 2 - Bitcoin
Set-Price 2 "34.086.91"
Set-Volume 2 "  -1.44%  "

# Row 3 - Ethereum
Set-Price 3 "1.793.15"
Set-Volume 3 "  -0.54%  "

# Row 4 - TetherUSD
Set-Volume 4 "  -0.04%  "

# Row 5 - BNB
Set-Price 5 "224.74"
Set-Volume 5 "  +0.28%  "

# Row 7 - USDC
Set-Volume 7 "  +0.00%  "

# Row 8 - Solana
Set-Price 8 "32.41"
Set-Volume 8 "  -0.07%  "

# Row 9 - Cardano
Set-Volume 9 "  -1.82%  "

# Row 10 - Dogecoin
Set-Volume 10 "  -0.42%  "

# Row 11 - TRON
Set-Price 11 "0.0929"
Set-Volume 11 "  +0.03%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-Price 12 "2.050.00"
Set-Volume 12 "  -0.60%  "

# Row 13 - WrappedEther
Set-Price 13 "1.782.82"
Set-Volume 13 "  -1.07%  "

# Row 14 - Chainlink
Set-Price 14 "10.85"
Set-Volume 14 "  -2.35%  "

# Row 15 - Polygon
Set-Price 15 "0.625"
Set-Volume 15 "  -2.58%  "

# Row 16 - WrappedBTC
Set-Price 16 "34.053.19"
Set-Volume 16 "  -1.65%  "

# Row 17 - Polkadot
Set-Volume 17 "  -3.58%  "

# Row 18 - Litecoin
Set-Price 18 "67.95"
Set-Volume 18 "  -1.79%  "

# Row 19 - BitcoinCash
Set-Price 19 "243.27"
Set-Volume 19 "  -3.55%  "

# Row 20 - ShibaInu
Set-Price 20 "0.0₃0784"
Set-Volume 20 "  -2.23%  "

# Row 21 - Dai
Set-Volume 21 "  +0.05%  "

# Row 22 - Avalanche
Set-Price 22 "10.69"
Set-Volume 22 "  -3.48%  "

# Row 23 - Uniswap
Set-Price 23 "4.08"
Set-Volume 23 "  -4.01%  "

# Row 25 - Monero
Set-Price 25 "158.99"
Set-Volume 25 "  -1.66%  "

# Row 26 - EthereumClassic
Set-Price 26 "16.25"
Set-Volume 26 "  -0.95%  "

# Row 27 - Cosmos
Set-Volume 27 "  -1.87%  "

# Row 28 - Stellar
Set-Volume 28 "  -1.96%  "

# Row 29 - BinanceUSD
Set-Volume 29 "  -0.01%  "

# Row 30 - Hedera
Set-Volume 30 "  -1.93%  "

# Row 31 - PancakeSwap
Set-Price 31 "1.22"
Set-Volume 31 "  +1.03%  "

# Row 32 - Filecoin
Set-Volume 32 "  -3.86%  "

# Row 33 - InternetComputer(DFINITY)
Set-Volume 33 "  -3.55%  "

# Row 34 - LidoDAOToken
Set-Volume 34 "  -3.90%  "

# Row 35 - Maker
Set-Price 35 "1.385.32"
Set-Volume 35 "  -3.26%  "

# Row 36 - ImmutableX
Set-Price 36 "0.645"
Set-Volume 36 "  +0.50%  "

# Row 37 - TrustWalletToken
Set-Volume 37 "  -1.99%  "

# Row 38 - VeChain
Set-Volume 38 "  -3.61%  "

# Row 39 / 40 - HuobiToken and Aave swap ranking positions
$ws.Cells.Item(39, 2).Value = "Aave"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-Price 39 "79.16"
Set-Volume 39 "  -6.53%  "

$ws.Cells.Item(40, 2).Value = "HuobiToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-Price 40 "2.35"
Set-Volume 40 "  -0.25%  "

# Row 41 - MXToken
Set-Volume 41 "  -3.57%  "

# Row 42 - ARBITRUM
Set-Price 42 "0.915"
Set-Volume 42 "  -4.80%  "

# Row 43 - RenderToken
Set-Price 43 "2.17"
Set-Volume 43 "  +0.35%  "

# Row 44 - BabyDogeCoin
Set-Price 44 "0.0₆0138"
Set-Volume 44 "  +9.00%  "

# Row 45 - Kaspa
Set-Price 45 "0.0495"
Set-Volume 45 "  -0.69%  "

# Row 46 - WEMIXToken
Set-Volume 46 "  -0.98%  "

# Row 47 - Quant
Set-Volume 47 "  +0.40%  "

# Row 48 - RocketPoolETH
Set-Price 48 "1.950.55"
Set-Volume 48 "  -0.29%  "

# Row 50 - PaxDollar
Set-Price 50 "0.999"
Set-Volume 50 "  -0.16%  "

# Row 51 - InjectiveProtocol
Set-Price 51 "11.94"
Set-Volume 51 "  -3.17%  "
